# Updated test data for German, Czech market; added test data for Belgium market.
#
#  - "Sounders" sheet is renamed to "Germany" (same data/layout as before).
#  - A brand-new "Belgium" sheet is inserted between "Germany" and "Czech",
#    built by duplicating the "Germany" sheet (so it picks up identical
#    column widths/styles/merged cells) and then updating its market name
#    and ticket reference.
#  - "Czech" is left alone content-wise; it simply shifts one slot to the
#    right in the tab order and stops being the active tab.

$wb = $excel.ActiveWorkbook

$germany = $wb.Worksheets.Item("Sounders")

# Duplicate "Sounders" right after itself -> becomes the new "Belgium" sheet,
# landing between Germany and Czech with the same formatting as Germany.
$germany.Copy($null, $germany) | Out-Null
$belgium = $wb.Worksheets.Item("Sounders (2)")

# Rename sheets.
$germany.Name = "Germany"
$belgium.Name = "Belgium"

# Fill in Belgium-specific market data (B2 = market name, B4 = ticket ref).
$belgium.Range("B2").Value = "Belgium Market"
$belgium.Range("B4").Value = "NGC-3478/T2269"

# Germany's selection becomes a full-sheet (select-all) selection.
$germany.Range("A1:XFD1048576").Select() | Out-Null

# Belgium is the newly active tab, with A4 selected.
$belgium.Activate()
$belgium.Range("A4").Select() | Out-Null
